$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 13 data (appended after existing data which ends at row 12)
$ws.Range("A13").Value = 8060.21
$ws.Range("B13").Value = 8019.31
$ws.Range("C13").Value = 17.8
$ws.Range("D13").Value = 17.89
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = 0.51
$ws.Range("G13").Value = 42620.766342592593
$ws.Range("H13").Value = $true

# Copy the date/time number format from the cell above (G12) so the new
# cell reuses the existing style instead of creating a new one.
$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
